$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I3").Value = "sv"
$ws.Range("J3").Value = "Statement-opinion"

# Row 11: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I11").Value = "sv"
$ws.Range("J11").Value = "Statement-opinion"

# Row 18: aa/Agree-Accept -> sd/Statement-non-opinion
$ws.Range("I18").Value = "sd"
$ws.Range("J18").Value = "Statement-non-opinion"

# Row 33: sd/Statement-non-opinion -> b/Acknowledge (Backchannel)
$ws.Range("I33").Value = "b"
$ws.Range("J33").Value = "Acknowledge (Backchannel)"

# Row 43: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I43").Value = "sv"
$ws.Range("J43").Value = "Statement-opinion"

$wb.Save()
